# actualizo lista de bugs
# Applies the BUGS.xlsx update: marks bug #6 (row 8, "Taggle") as CORREGIDO
# (previously PENDIENTE), and adds a new bug #15 in row 17 about the
# "fases" object returning groups with all-identical teams.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (bug #1): status flips from PENDIENTE to CORREGIDO -------------
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F3").Value = "CORREGIDO"

# --- Row 8 (bug #6, "Taggle"): status flips from PENDIENTE to CORREGIDO ---
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F8").Value = "CORREGIDO"

# --- New bug row 17 (id 15) -------------------------------------------------
$ws.Rows.Item(17).RowHeight = 42

$ws.Range("A17").Value = 15

$ws.Range("B17").Value = "Los grupos se guardan con los mismos equipos"
$ws.Range("B17").Font.Color = 255
$ws.Range("B17").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B17").VerticalAlignment = -4108     # xlCenter
$ws.Range("B17").WrapText = $true

$ws.Range("E17").Value = "edicion-fases.aspx"

$ws.Range("C17").Value = "Cuando crea los grupos de equipos, el objeto fases, devuelve grupos con todos equipos iguales :( "
$ws.Range("C17").WrapText = $true

$ws.Range("D17").Value = "Tony"

$ws.Range("F4").Copy() | Out-Null
$ws.Range("F17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F17").Value = "PENDIENTE"

# --- View / selection ------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H8").Select() | Out-Null
